# The commit adds eight new "game details" lines (Developer, Publisher,
# Genres, Release year, Critics score, Min. CPU, Min. videocard, Min. RAM)
# right after the existing "Название игры: <GAME_NAME>" paragraph, followed
# by one trailing blank paragraph, all styled like the rest of the document
# (Times New Roman, sz 28) with per-run language tagging (Russian labels,
# English placeholders).

$d = $word.ActiveDocument

# --- locate the paragraph holding the <GAME_NAME> placeholder -------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*<GAME_NAME>*") {
        $targetPara = $p
        break
    }
}

# --- small helpers to build well-formed WordprocessingML fragments --------
function New-Run {
    param(
        [string]$Lang,
        [string]$Text,
        [bool]$Preserve = $false
    )
    $space = ""
    if ($Preserve) { $space = ' xml:space="preserve"' }
    $rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="' + $Lang + '"/></w:rPr>'
    return '<w:r>' + $rPr + '<w:t' + $space + '>' + $Text + '</w:t></w:r>'
}

function New-Para {
    param(
        [string]$RunsXml
    )
    $pPr = '<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr>'
    return '<w:p>' + $pPr + $RunsXml + '</w:p>'
}

# --- the eight content paragraphs ------------------------------------------
$paraDeveloper = New-Para ((New-Run "ru-RU" "Разработчик:") + (New-Run "en-US" " &lt;DEVELOPER&gt;" $true))

$paraPublisher = New-Para ((New-Run "ru-RU" "Публикатор:") + (New-Run "en-US" " &lt;PUBLISHER&gt;" $true))

$paraGenres = New-Para ((New-Run "ru-RU" "Жанры") + (New-Run "en-US" ": &lt;GENRES&gt;"))

$paraReleaseYear = New-Para ((New-Run "ru-RU" "Год выхода:") + (New-Run "en-US" " &lt;RELEASE_YEAR&gt;" $true))

$paraCriticsScore = New-Para ((New-Run "ru-RU" "Оценка критиков") + (New-Run "en-US" ":") + (New-Run "en-US" " &lt;CRITICS_SCORE&gt;" $true))

$paraCpu = New-Para ((New-Run "ru-RU" "Мин. " $true) + (New-Run "ru-RU" "п") + (New-Run "ru-RU" "роцессор:") + (New-Run "en-US" " &lt;CPU_NAME&gt;" $true))

$paraVideocard = New-Para ((New-Run "ru-RU" "Мин. " $true) + (New-Run "ru-RU" "видеокарта") + (New-Run "ru-RU" ":") + (New-Run "en-US" " &lt;VIDEOCARD_NAME&gt;" $true))

$paraRam = New-Para ((New-Run "ru-RU" "Количество оперативной памяти") + (New-Run "en-US" ":") + (New-Run "en-US" " &lt;RAM_AMOUNT&gt;" $true))

# trailing empty paragraph (must be the *last* element of the fragment: Word
# merges the very last inserted paragraph mark into whatever follows the
# insertion point, so putting the empty paragraph last makes it come out as
# its own standalone paragraph instead of swallowing "Цена: <PRICE>").
$paraTrailingBlank = New-Para ""

$bodyXml = $paraDeveloper + $paraPublisher + $paraGenres + $paraReleaseYear + $paraCriticsScore + $paraCpu + $paraVideocard + $paraRam + $paraTrailingBlank

$xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- insert right after the <GAME_NAME> paragraph --------------------------
$insertPos = $targetPara.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertXML($xmlFragment)
